$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Year of Treatment" column (B); this shifts C:H left into B:G,
# so the former "Total" column (H) becomes the new last column (G).
$ws.Columns.Item(2).Delete()

# Retitle the (now shifted) header cells.
$ws.Range("B1").Value = "Never went to school / never completed primary school (ISCED 0).deja.deja.deja"
$ws.Range("C1").Value = "Primary level of education (ISCED 1).deja.deja.deja"
$ws.Range("D1").Value = "Secondary level of education (ISCED 2 and ISCED 3).deja.deja.deja"
$ws.Range("E1").Value = "Higher education (ISCED 4 to ISCED 6).deja.deja.deja"
$ws.Range("F1").Value = "Not known / missing.deja.deja.deja"
$ws.Range("G1").Value = "Total.deja.deja.deja"
